$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert 4 new blank rows right before the current row 221 (A1001-style
# key/value block), shifting the existing rows 221+ down to 225+.
$ws.Range("A218:A221").EntireRow.Insert()

# Fill the newly inserted rows with the new key/value pairs (prob x28).
# Column A first (in row order), then column B (in row order), then the
# final key cell -- this mirrors the order the strings were authored in
# (and therefore the order they land in the shared string table).
$ws.Range("A218").Value2 = "x1008"
$ws.Range("A219").Value2 = "x1009"
$ws.Range("A220").Value2 = "x1010"

$ws.Range("B218").Value2 = "함수의 정의역의 각 원소에 대한 함숫값의 범위를 조사합니다."
$ws.Range("B219").Value2 = "함수의 조건에 맞는 가능한 치역을 모두 구합니다."
$ws.Range("B220").Value2 = "각 치역에 대해 조건에 맞는 함수의 개수를 구합니다."
$ws.Range("B221").Value2 = "각각의 개수를 모두 더해서 조건을 만족시키는 전체 개수를 구합니다."

$ws.Range("A221").Value2 = "x1011"

# Restore the view: scrolled so row 208 is at the top, with C221 selected.
$ws.Activate()
$ws.Range("C221").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 208
$win.ScrollColumn = 1
